$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = ""

$ws.Range("H19").Value = 1428.5714
$ws.Range("I19").Value = 750
$ws.Range("J19").Value = 1700
$ws.Range("K19").Value = 750
$ws.Range("L19").Value = 1700
$ws.Range("M19").Value = -575
$ws.Range("N19").Value = -2050

$ws.Range("H98").Value = 37912840
$ws.Range("I98").Value = 40281896
$ws.Range("J98").Value = 8000
$ws.Range("K98").Value = 40281896
$ws.Range("L98").Value = 8000
$ws.Range("M98").Value = -40280398
$ws.Range("N98").Value = -10996

$ws.Range("H100").Value = 19610296
$ws.Range("I100").Value = 33335384
$ws.Range("J100").Value = 3029.4285
$ws.Range("K100").Value = 33335384
$ws.Range("L100").Value = 3029.4285
$ws.Range("M100").Value = -33334843
$ws.Range("N100").Value = -4111.4285

$ws.Range("H111").Value = 848.5294
$ws.Range("I111").Value = 836.1111
$ws.Range("J111").Value = 862.5
$ws.Range("K111").Value = 2508.3333
$ws.Range("L111").Value = 2587.5
$ws.Range("M111").Value = 558.6667000000002
$ws.Range("N111").Value = -8721.5

$ws.Range("H122").Value = 37912840
$ws.Range("I122").Value = 40281896
$ws.Range("J122").Value = 8000
$ws.Range("K122").Value = 120845688
$ws.Range("L122").Value = 24000
$ws.Range("M122").Value = -120843238
$ws.Range("N122").Value = -28900

$ws.Range("H132").Value = 2526024.8
$ws.Range("I132").Value = 2598182.8
$ws.Range("J132").Value = 500
$ws.Range("K132").Value = 7794548.399999999
$ws.Range("L132").Value = 1500
$ws.Range("M132").Value = -7792018.399999999
$ws.Range("N132").Value = -6560

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H49").Value = 5000
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 5000
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 5000
$ws.Range("N49").Value = -5520

$ws.Range("H61").Value = 1658.0769
$ws.Range("I61").Value = 1427.8667
$ws.Range("J61").Value = 1972
$ws.Range("K61").Value = 1427.8667
$ws.Range("L61").Value = 1972
$ws.Range("M61").Value = -1215.8667
$ws.Range("N61").Value = -2396

$ws.Range("H96").Value = 28000
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 28000
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 28000
$ws.Range("N96").Value = -33492

$ws.Range("H136").Value = 1658.0769
$ws.Range("I136").Value = 1427.8667
$ws.Range("J136").Value = 1972
$ws.Range("K136").Value = 4283.6001
$ws.Range("L136").Value = 5916
$ws.Range("M136").Value = -1733.6001
$ws.Range("N136").Value = -11016

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H95").Value = 27777
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 27777
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 27777
$ws.Range("N95").Value = -33269

$ws.Range("H134").Value = 20925.555
$ws.Range("I134").Value = 5788
$ws.Range("J134").Value = 73907
$ws.Range("K134").Value = 17364
$ws.Range("L134").Value = 221721
$ws.Range("M134").Value = -14829
$ws.Range("N134").Value = -226791

$ws.Range("H137").Value = 59800
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 59800
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 59800
$ws.Range("N137").Value = -70000

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1768478.4
$ws.Range("I31").Value = 2464747.2
$ws.Range("J31").Value = 85828.5
$ws.Range("K31").Value = 2464747.2
$ws.Range("L31").Value = 85828.5
$ws.Range("M31").Value = -2464452.2

$ws.Range("H34").Value = 1768478.4
$ws.Range("I34").Value = 2464747.2
$ws.Range("J34").Value = 85828.5
$ws.Range("K34").Value = 2464747.2
$ws.Range("L34").Value = 85828.5
$ws.Range("M34").Value = -2464545.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1021
$ws.Range("I5").Value = 829.6111
$ws.Range("J5").Value = 1710
$ws.Range("K5").Value = 2488.8333
$ws.Range("L5").Value = 5130
$ws.Range("M5").Value = -2376.8333
$ws.Range("N5").Value = -5354

$ws.Range("H122").Value = 570.5263
$ws.Range("I122").Value = 332.66666
$ws.Range("J122").Value = 784.6
$ws.Range("K122").Value = 2993.99994
$ws.Range("L122").Value = 7061.400000000001
$ws.Range("M122").Value = -543.9999399999997
$ws.Range("N122").Value = -11961.4

$ws.Range("H126").Value = 2678.75
$ws.Range("I126").Value = 965
$ws.Range("J126").Value = 3250
$ws.Range("K126").Value = 2895
$ws.Range("L126").Value = 9750
$ws.Range("M126").Value = 2045
$ws.Range("N126").Value = -19630

$ws.Range("H135").Value = 1021
$ws.Range("I135").Value = 829.6111
$ws.Range("J135").Value = 1710
$ws.Range("K135").Value = 7466.4999
$ws.Range("L135").Value = 15390
$ws.Range("M135").Value = -4931.4999
$ws.Range("N135").Value = -20460

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1483.1111
$ws.Range("I113").Value = 1409.8
$ws.Range("J113").Value = 1574.75
$ws.Range("K113").Value = 1409.8
$ws.Range("L113").Value = 1574.75
$ws.Range("M113").Value = 760.2
$ws.Range("N113").Value = -5914.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H34").Value = 3861.125
$ws.Range("I34").Value = 3450
$ws.Range("J34").Value = 3998.1667
$ws.Range("K34").Value = 3450
$ws.Range("L34").Value = 3998.1667
$ws.Range("M34").Value = -3278
$ws.Range("N34").Value = -4342.1667

$ws.Range("H39").Value = 4160
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 4160
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 4160
$ws.Range("M39").Value = ""
$ws.Range("N39").Value = -5080

$ws.Range("H40").Value = 1872.7878
$ws.Range("I40").Value = 1753.88
$ws.Range("J40").Value = 2244.375
$ws.Range("K40").Value = 1753.88
$ws.Range("L40").Value = 2244.375
$ws.Range("M40").Value = -1617.88
$ws.Range("N40").Value = -2516.375

$ws.Range("H46").Value = 6768.294
$ws.Range("I46").Value = 892.38464
$ws.Range("J46").Value = 25865
$ws.Range("K46").Value = 892.38464
$ws.Range("L46").Value = 25865
$ws.Range("M46").Value = -704.38464
$ws.Range("N46").Value = -26241

$ws.Range("H61").Value = 3000
$ws.Range("I61").Value = 2000
$ws.Range("J61").Value = 3333.3333
$ws.Range("K61").Value = 2000
$ws.Range("L61").Value = 3333.3333
$ws.Range("M61").Value = -1798
$ws.Range("N61").Value = -3737.3333

$ws.Range("H113").Value = 3000
$ws.Range("I113").Value = 2000
$ws.Range("J113").Value = 3333.3333
$ws.Range("K113").Value = 2000
$ws.Range("L113").Value = 3333.3333
$ws.Range("M113").Value = 170
$ws.Range("N113").Value = -7673.3333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 59621.41
$ws.Range("I136").Value = 67437.6
$ws.Range("J136").Value = 1000
$ws.Range("K136").Value = 202312.8
$ws.Range("L136").Value = 3000
$ws.Range("M136").Value = -199762.8
$ws.Range("N136").Value = -8100
